$wb = $excel.ActiveWorkbook

# --- Update the Date value on the Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-04-11T11:42:22-04:00"

# --- Copy Display (column C) into Definition (column D) on the Concepts sheet ---
$concepts = $wb.Worksheets.Item("Concepts")
for ($row = 2; $row -le 12; $row++) {
    $display = $concepts.Cells.Item($row, 3).Text
    $concepts.Cells.Item($row, 4).Value = $display
}
